# FLDEP 20210408 data dump
# 1) Re-sort the existing station table (rows 2-20) alphabetically by Name
#    (column C), which is what Excel's "Sort" command does and is how the
#    rows ended up reshuffled in the saved file.
# 2) Append the newly-dumped "Piney 17-23" stations as rows 21-27.
# 3) Leave the final selection on the last entered cell, as in the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: sort A2:E20 ascending by column C (Name) ---
$sortRange = $ws.Range("A2:E20")
$sortKey   = $ws.Range("C2:C20")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortKey, 0, 1, 0, 0) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# --- Step 2: append the new data-dump rows ---
$newRows = @(
    @(27.781299000000001, -82.474100000000007, "Piney 22"),
    @(27.682528999999999, -82.496769999999998, "Piney 18"),
    @(27.728999999999999, -82.498699999999999, "Piney 20"),
    @(27.723801000000002, -82.533799999999999, "Piney 19"),
    @(27.777999999999999, -82.520300000000006, "Piney 21"),
    @(27.776399999999999, -82.438730000000007, "Piney 23"),
    @(27.693398999999999, -82.555899999999994, "Piney 17")
)

$row = 21
foreach ($entry in $newRows) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

# --- Step 3: match the saved file's final selection ---
$ws.Range("E26").Select()
